$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 256 (shifts existing rows 256:286 down to 257:287)
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row with the new observation
$ws.Range("A256").Value = 2
$ws.Range("B256").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C256").Value = "Coquimbo"
$ws.Range("D256").Value = 45077
$ws.Range("E256").Value = 4
$ws.Range("F256").Value = 100112031
$ws.Range("G256").Value = "Poroto verde"
$ws.Range("H256").Value = "Magnum"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 1100
$ws.Range("K256").Value = 14000
$ws.Range("L256").Value = 16000
$ws.Range("M256").Value = 15000
$ws.Range("N256").Value = "$/malla 25 kilos"
$ws.Range("O256").Value = "Provincia de Limarí"
$ws.Range("P256").Value = 600
$ws.Range("Q256").Value = 25
$ws.Range("R256").Value = "Hortaliza"
